$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.939.42'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.886.85'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.019'
$ws.Range('E4').Value = '  +1.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.53'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.019'
$ws.Range('E6').Value = '  +1.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4650'
$ws.Range('E7').Value = '  -1.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3905'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.98'
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07960'
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.59'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.907.25'
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.082'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.019'
$ws.Range('E16').Value = '  +1.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06772'
$ws.Range('E17').Value = '  +2.33%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '86.91'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001044'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.05'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('E21').Value = '  +1.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.966.32'
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.477'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.352'
$ws.Range('E25').Value = '  +2.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.122.17'
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.17'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.94'
$ws.Range('E28').Value = '  -1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.062'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.401'
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '121.54'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9584'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09478'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.674'
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.320'
$ws.Range('E35').Value = '  +0.51%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.356'
$ws.Range('E36').Value = '  -6.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06095'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02233'
$ws.Range('E38').Value = '  -0.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.215'
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.067'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5934'
$ws.Range('E41').Value = '  -0.88%  '
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.23'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.273'
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5650'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.11'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.394'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06902'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '113.56'
$ws.Range('E50').Value = '  +2.35%  '
$ws.Range('E51').Value = '  -0.29%  '
